$d = $word.ActiveDocument
Write-Output "Frames.Count: $($d.Frames.Count)"
